$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the header values in F1:J1 (keep cell formatting/style)
$ws.Range("F1:J1").ClearContents()

# Update the active selection to L6
$ws.Range("L6").Select()
